$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 257, shifting the existing weekly records
# (previously rows 257-272) down to rows 258-273.
$ws.Rows(257).Insert()

# Populate the new row with the latest weekly price record for
# "Ajo" (Garlic) - Chino, Primera quality - matching the constant
# columns used throughout this data block.
$ws.Cells.Item(257, 1).Value = 7
$ws.Cells.Item(257, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(257, 3).Value = "Ñuble"
$ws.Cells.Item(257, 4).Value = 44826
$ws.Cells.Item(257, 5).Value = 16
$ws.Cells.Item(257, 6).Value = 100112003
$ws.Cells.Item(257, 7).Value = "Ajo"
$ws.Cells.Item(257, 8).Value = "Chino"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 60
$ws.Cells.Item(257, 11).Value = 23000
$ws.Cells.Item(257, 12).Value = 24000
$ws.Cells.Item(257, 13).Value = 23500
$ws.Cells.Item(257, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(257, 15).Value = "China"
$ws.Cells.Item(257, 16).Value = 2350
$ws.Cells.Item(257, 17).Value = 10
$ws.Cells.Item(257, 18).Value = "Hortaliza"
